$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to be interpreted/stored as text, matching the
    # inlineStr cells in the source workbook, then strip the temporary
    # text number-format so no stray style gets attached to the cell.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" "61.165.39"
$ws.Range("E2").Value = "  +0.05%  "

Set-TextValue "D3" "3.374.90"
$ws.Range("E3").Value = "  +1.88%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue "D5" "570.97"
$ws.Range("E5").Value = "  +0.21%  "

Set-TextValue "D6" "137.40"
$ws.Range("E6").Value = "  +8.31%  "

$ws.Range("E7").Value = "  +0.02%  "

Set-TextValue "D8" "3.375.72"
$ws.Range("E8").Value = "  +1.90%  "

Set-TextValue "D9" "0.474"
$ws.Range("E9").Value = "  -0.31%  "

$ws.Range("E10").Value = "  +4.01%  "

$ws.Range("E11").Value = "  +5.20%  "

$ws.Range("E12").Value = "  +5.15%  "

Set-TextValue "D13" "3.956.69"
$ws.Range("E13").Value = "  +2.01%  "

$ws.Range("E15").Value = "  +3.40%  "

Set-TextValue "D16" "3.384.05"
$ws.Range("E16").Value = "  +2.15%  "

Set-TextValue "D17" "25.16"
$ws.Range("E17").Value = "  +2.84%  "

Set-TextValue "D18" "61.288.75"
$ws.Range("E18").Value = "  +0.13%  "

Set-TextValue "D19" "13.93"
$ws.Range("E19").Value = "  +5.92%  "

Set-TextValue "D20" "5.85"
$ws.Range("E20").Value = "  +4.78%  "

$ws.Range("E21").Value = "  +4.38%  "

Set-TextValue "D22" "381.58"
$ws.Range("E22").Value = "  +8.60%  "

Set-TextValue "D23" "0.574"
$ws.Range("E23").Value = "  +4.13%  "

Set-TextValue "D24" "3.512.82"
$ws.Range("E24").Value = "  +1.94%  "

Set-TextValue "D25" "1.00"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  +1.39%  "

$ws.Range("E27").Value = "  +11.81%  "

Set-TextValue "D28" "1.67"
$ws.Range("E28").Value = "  +15.97%  "

$ws.Range("E29").Value = "  +8.44%  "

$ws.Range("E30").Value = "  +0.26%  "

Set-TextValue "D31" "8.16"
$ws.Range("E31").Value = "  +4.41%  "

$ws.Range("E32").Value = "  +5.69%  "

Set-TextValue "D33" "2.13"
$ws.Range("E33").Value = "  +1.35%  "

Set-TextValue "D35" "3.407.39"
$ws.Range("E35").Value = "  +1.89%  "

Set-TextValue "D36" "23.43"
$ws.Range("E36").Value = "  +4.75%  "

$ws.Range("E37").Value = "  +4.73%  "

Set-TextValue "D38" "7.01"
$ws.Range("E38").Value = "  +5.49%  "

$ws.Range("E39").Value = "  +5.36%  "

Set-TextValue "D40" "162.59"
$ws.Range("E40").Value = "  +0.19%  "

Set-TextValue "D41" "0.0804"
$ws.Range("E41").Value = "  +7.41%  "

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("E43").Value = "  +4.96%  "

Set-TextValue "D44" "41.49"
$ws.Range("E44").Value = "  +1.17%  "

Set-TextValue "D45" "0.759"
$ws.Range("E45").Value = "  +2.24%  "

$ws.Range("E46").Value = "  +7.92%  "

Set-TextValue "D47" "1.69"

Set-TextValue "D48" "23.14"
$ws.Range("E48").Value = "  +4.82%  "

Set-TextValue "D49" "6.94"
$ws.Range("E49").Value = "  +4.67%  "

Set-TextValue "D50" "23.10"
$ws.Range("E50").Value = "  +11.63%  "

Set-TextValue "D51" "2.338.02"
$ws.Range("E51").Value = "  +6.30%  "
